$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set up the "year" label cells (column A) by copying the existing
# label cell's formatting (bold, centered, bordered) from row 10, then
# overwrite with the new year text. ---
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("A11").Value = "2021年"
$ws.Range("A12").Value = "2022年"

# --- Row 11: 2021年 ---
$ws.Range("B11").Value = 221.48
$ws.Range("C11").Value = 106.59
$ws.Range("D11").Value = 51.96
$ws.Range("F11").Value = 105.1
$ws.Range("G11").Value = 2302.12
$ws.Range("H11").Value = 105
$ws.Range("I11").Value = 1407.58
$ws.Range("J11").Value = 52.65
$ws.Range("K11").Value = 24435.19
$ws.Range("L11").Value = 27.5
$ws.Range("M11").Value = 49.84
$ws.Range("N11").Value = -57.68
$ws.Range("O11").Value = 35.41
$ws.Range("P11").Value = 1062.46
$ws.Range("Q11").Value = 252.63
$ws.Range("R11").Value = 2.28
$ws.Range("S11").Value = 50.14
$ws.Range("T11").Value = 242.04
$ws.Range("U11").Value = 2553.64
$ws.Range("V11").Value = 1182.03
$ws.Range("W11").Value = 4267.56
$ws.Range("X11").Value = 435.54
$ws.Range("Y11").Value = 1206.72
$ws.Range("Z11").Value = 198.46
$ws.Range("AA11").Value = 3.35
$ws.Range("AB11").Value = 1437.91
$ws.Range("AC11").Value = 682.24
$ws.Range("AD11").Value = 26.21
$ws.Range("AE11").Value = 9.470000000000001
$ws.Range("AF11").Value = 1007.39
$ws.Range("AG11").Value = 331.51
$ws.Range("AH11").Value = 55.45
$ws.Range("AI11").Value = 1404.01
$ws.Range("AJ11").Value = 14.73
$ws.Range("AK11").Value = 213.92
$ws.Range("AL11").Value = 197.25
$ws.Range("AM11").Value = 1172.51
$ws.Range("AN11").Value = 137.76
$ws.Range("AO11").Value = 55.13
$ws.Range("AP11").Value = 1459.97
$ws.Range("AQ11").Value = 365.06

# E11 has no reported value for this indicator/year; leave it blank
# (mirrors the blank-cell representation already used for E10).
$ws.Range("E10").Copy($ws.Range("E11"))

# --- Row 12: 2022年 -- only the overall total (column K) has been
# reported so far; every other indicator cell stays blank. ---
$ws.Range("E10").Copy($ws.Range("B12:AQ12"))
$ws.Range("K12").Value = 23792
